$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 72
$prev = $row - 1

# Copy formats (styles) from the previous row so A72 gets style 1 and E72 gets style 2
$ws.Range("A$prev`:V$prev").Copy()
$ws.Range("A$row").PasteSpecial(-4122)

$ws.Cells.Item($row, 1).Value = 71
$ws.Cells.Item($row, 2).Value = "armenia"
$ws.Cells.Item($row, 3).Value = "premier-league"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45234.5
$ws.Cells.Item($row, 6).Value = "Alashkert"
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = "Shirak Gyumri"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 1.4
$ws.Cells.Item($row, 11).Value = "03/11/2023 00:12"
$ws.Cells.Item($row, 12).Value = 1.38
$ws.Cells.Item($row, 13).Value = "04/11/2023 11:50"
$ws.Cells.Item($row, 14).Value = 4.38
$ws.Cells.Item($row, 15).Value = "03/11/2023 00:12"
$ws.Cells.Item($row, 16).Value = 4.73
$ws.Cells.Item($row, 17).Value = "04/11/2023 11:50"
$ws.Cells.Item($row, 18).Value = 6.17
$ws.Cells.Item($row, 19).Value = "03/11/2023 00:12"
$ws.Cells.Item($row, 20).Value = 8.26
$ws.Cells.Item($row, 21).Value = "04/11/2023 11:50"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/armenia/premier-league/alashkert-shirak-gyumri/0v7InesP/"
